$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.418.76'
$ws.Range("E2").Value = '  +1.95%  '
$ws.Range("D3").Value = '2.596.46'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.96'
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.12'
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("E8").Value = '  +0.42%  '
$ws.Range("D9").Value = '2.616.93'
$ws.Range("E9").Value = '  +1.08%  '
$ws.Range("E10").Value = '  -0.77%  '
$ws.Range("E11").Value = '  +1.88%  '
$ws.Range("E12").Value = '  +4.18%  '
$ws.Range("E13").Value = '  -6.62%  '
$ws.Range("D14").Value = '3.057.84'
$ws.Range("E14").Value = '  +0.73%  '
$ws.Range("D15").Value = '60.411.59'
$ws.Range("E15").Value = '  +2.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.39'
$ws.Range("E16").Value = '  +0.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000140'
$ws.Range("E17").Value = '  +2.79%  '
$ws.Range("D18").Value = '2.608.05'
$ws.Range("E18").Value = '  +0.94%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.29'
$ws.Range("E19").Value = '  +9.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.66'
$ws.Range("E20").Value = '  +1.78%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '346.71'
$ws.Range("E21").Value = '  +2.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.98'
$ws.Range("E22").Value = '  +9.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.532'
$ws.Range("E24").Value = '  +14.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.16'
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("E27").Value = '  -1.85%  '
$ws.Range("E28").Value = '  +4.86%  '
$ws.Range("E29").Value = '  +1.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.82'
$ws.Range("E30").Value = '  +9.21%  '
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("E32").Value = '  +3.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '160.86'
$ws.Range("E33").Value = '  +0.34%  '
$ws.Range("E34").Value = '  +2.59%  '
$ws.Range("E35").Value = '  +4.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.960'
$ws.Range("E36").Value = '  +9.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.22'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.61'
$ws.Range("E38").Value = '  +8.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.78'
$ws.Range("E39").Value = '  +1.06%  '
$ws.Range("E40").Value = '  +3.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.848'
$ws.Range("E41").Value = '  -2.62%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '137.76'
$ws.Range("E43").Value = '  +4.19%  '
$ws.Range("E44").Value = '  -0.27%  '
$ws.Range("E45").Value = '  +2.06%  '
$ws.Range("E46").Value = '  +1.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.65'
$ws.Range("E47").Value = '  +3.52%  '
$ws.Range("E48").Value = '  +1.93%  '
$ws.Range("E49").Value = '  +3.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.86'
$ws.Range("E50").Value = '  +6.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.68'
$ws.Range("E51").Value = '  +0.27%  '
